$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
# Columns A-J (1-10): "<name>_old"  -> "<name>_FV2210"
$oldHeaders = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

# Column K (11) "diff" is unchanged.

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2304"
$newHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Turn the data range into an Excel Table (ListObject) with autofilter ---
$range = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split/freeze at row 2) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
